# Apply futures-roll / side-flip updates to the trade file (PROD rec and logs)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 : ADU5 -> ADZ5 (Curncy ticker roll) ---
$ws.Range("C2").Value = "ADZ5 Curncy"

# --- Row 4 : BPU5 -> BPZ5, side flip SELL -> BUY, target -14 -> 14 ---
$ws.Range("C4").Value = "BPZ5 Curncy"
$ws.Range("F4").Value = "BUY"
$ws.Range("O4").Value = 14

# --- Row 5 : side flip BUY -> SELL, target 16 -> -16 ---
$ws.Range("F5").Value = "SELL"
$ws.Range("O5").Value = -16

# --- Row 6 : side flip BUY -> SELL, target 10 -> -10 ---
$ws.Range("F6").Value = "SELL"
$ws.Range("O6").Value = -10

# --- Row 7 : ESU5 -> ESZ5 (Index ticker roll) ---
$ws.Range("C7").Value = "ESZ5 Index"

# --- Row 8 : ECU5 -> ECZ5 (Curncy ticker roll) ---
$ws.Range("C8").Value = "ECZ5 Curncy"

# --- Row 10 : JYU5 -> JYZ5, side flip SELL -> BUY, target -14 -> 14 ---
$ws.Range("C10").Value = "JYZ5 Curncy"
$ws.Range("F10").Value = "BUY"
$ws.Range("O10").Value = 14

# --- Row 11 : NQU5 -> NQZ5, side flip SELL -> BUY, target -6 -> 6 ---
$ws.Range("C11").Value = "NQZ5 Index"
$ws.Range("F11").Value = "BUY"
$ws.Range("O11").Value = 6

# --- Row 12 : RTYU5 -> RTYZ5, quantity 27 -> 26, side flip SELL -> BUY, target -27 -> 26 ---
$ws.Range("C12").Value = "RTYZ5 Index"
$ws.Range("D12").Value = 26
$ws.Range("F12").Value = "BUY"
$ws.Range("O12").Value = 26

# --- Row 16 : quantity 33 -> 34, target 33 -> 34 ---
$ws.Range("D16").Value = 34
$ws.Range("O16").Value = 34

# --- Row 18 : RXU5 -> RXZ5 (Comdty ticker roll) ---
$ws.Range("C18").Value = "RXZ5 Comdty"

# --- Row 19 : OEU5 -> OZU5 (Comdty ticker roll) ---
$ws.Range("C19").Value = "OZU5 Comdty"

# --- Row 20 : CLV5 -> CLZ5 (Comdty ticker roll) ---
$ws.Range("C20").Value = "CLZ5 Comdty"

# --- Row 21 : side flip SELL -> BUY, target -3 -> 3 ---
$ws.Range("F21").Value = "BUY"
$ws.Range("O21").Value = 3

# --- Row 22 : side flip BUY -> SELL, target 10 -> -10 ---
$ws.Range("F22").Value = "SELL"
$ws.Range("O22").Value = -10

# --- Row 23 : NGV25 -> NGX25 (Comdty ticker roll) ---
$ws.Range("C23").Value = "NGX25 Comdty"

# --- Row 25 : side flip SELL -> BUY, target -5 -> 5 ---
$ws.Range("F25").Value = "BUY"
$ws.Range("O25").Value = 5
